$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "Khaleel Ahmed"
$ws.Name = "Khaleel Ahmed"

# Insert a new column before column A, shifting existing columns (A:L) to (B:M)
$ws.Columns("A:A").Insert()

# Fill in the new matchNo column
$ws.Range("A1").Value = "matchNo"
$ws.Range("A2").Value = "9th"
